# Insert a new column before column F ("eta") — this shifts the existing
# F:Q data (pt_min..diff) over to G:R and automatically rewrites every
# formula reference that pointed into the shifted range.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("F:F").Insert()

# Populate the new "eta" column (header + six data rows).
$ws.Range("F1").Value = "eta"
$ws.Range("F2").Value = -1.27
$ws.Range("F3").Value = -0.74
$ws.Range("F4").Value = -0.27
$ws.Range("F5").Value = 0.27
$ws.Range("F6").Value = 0.74
$ws.Range("F7").Value = 1.27

# The column insert leaves the third shared-formula chain (now in column R)
# expanded into one independent formula per row instead of a single shared
# group; re-apply the formula across the block so Excel re-collapses it
# back into a shared formula the same way it does for the other two chains.
$ws.Range("R3:R7").Formula = "=SQRT(O3*O3+P3*P3)"

# Match the author's final selection.
$null = $ws.Range("G23").Select()
